$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '255.13'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '3.96%'

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '27.60'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '-7.55%'

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '5.185'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '0.31%'

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.05862'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '2.05%'

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '6.717'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '1.08%'

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.8679'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '1.13%'

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.9538'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '11.76%'

$ws.Range("B9").Value = 'WazirX'
$ws.Range("C9").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.1409'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '1.98%'

$ws.Range("B10").Value = 'MandalaExchangeToken'
$ws.Range("C10").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07171'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '1.29%'

$ws.Range("B11").Value = 'BitrueCoin'
$ws.Range("C11").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.03179'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '1.04%'

$ws.Range("B12").Value = 'BitMartToken'
$ws.Range("C12").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.09233'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '-1.28%'

$ws.Range("B13").Value = 'BitForexToken'
$ws.Range("C13").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.001543'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '0.39%'

$ws.Range("B14").Value = 'One'
$ws.Range("C14").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0006080'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '-94.04%'

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.005969'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '-0.92%'

$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '-0.84%'

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.227'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '-1.78%'

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.205'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '0.62%'

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.03436'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '3.47%'

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.1281'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '-1.34%'

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '3.554'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '1.94%'

$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '1.52%'

$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '-2.48%'

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.001227'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '0.44%'

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.004790'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '15.18%'

$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '-0.01%'

$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '1.21%'

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.03809'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '1.64%'

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.005629'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '5.03%'

$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '3.09%'

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.002338'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '6.28%'

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.009798'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '-1.51%'

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.00005363'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '-1.60%'

$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '-0.01%'

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.1091'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '21.38%'

$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '-3.63%'

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.00002100'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '-0.01%'

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0002000'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '-0.01%'
